# Apply edits described by the diff:
#  - Insert a new worksheet "Player Info" before "ODI Batting"
#  - Populate "Player Info" with ID / NAME / BATTING_HAND / BOWL_STYLE data
#  - Update "ODI Batting" column D (MATCH_CARD_LINK -> MATCH_CODE), replacing
#    the scorecard URLs with their bare match-code values

$wb = $excel.ActiveWorkbook

# Existing sheet (the one currently named "ODI Batting") - make it active so
# the newly-added sheet is inserted immediately before it.
$wb.Worksheets.Item("ODI Batting").Activate()

# Add a brand-new worksheet; Excel inserts it immediately before the active sheet
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

# Re-fetch "ODI Batting" by name now that the sheet collection has shifted -
# earlier references are positional and go stale once a sheet is inserted.
$batting = $wb.Worksheets.Item("ODI Batting")

# --- Populate "Player Info" ---
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Match the bold/centered/bordered header formatting used on "ODI Batting" -
# copy it straight from that sheet's own header cell so we reuse the same
# style definition instead of inventing a new (merely equivalent) one.
$batting.Range("A1").Copy() | Out-Null
$playerInfo.Range("A1:D1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# "ID" is numeric-looking text in the source data ("6465", not 6465), so
# force text entry (NumberFormat "@") then drop the formatting again so the
# cell ends up plain/unstyled - same way the rest of this workbook's
# numeric-looking columns (MATCH_NUMBER, RUNS_SCORED, ...) are stored.
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "6465"
$playerInfo.Range("A2").ClearFormats()

$playerInfo.Range("B2").Value = "Justin Pierre Greaves"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium"

# --- Update "ODI Batting" column D: header rename + URLs -> bare codes ---
$batting.Range("D1").Value = "MATCH_CODE"

$batting.Range("D2").NumberFormat = "@"
$batting.Range("D2").Value = "4519"
$batting.Range("D2").ClearFormats()

$batting.Range("D3").NumberFormat = "@"
$batting.Range("D3").Value = "4520"
$batting.Range("D3").ClearFormats()

$batting.Range("D4").NumberFormat = "@"
$batting.Range("D4").Value = "4522"
$batting.Range("D4").ClearFormats()
